# Generate Report for Handback
# Adds two new handed-back files (2b431392-... and 310c7c43-...) to the
# Overview / zh-cn / de-de worksheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# Sheet "Overview": two new rows (6 and 7), columns A (file + hyperlink),
# B (zh-cn status), C (de-de status)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A6").Value = "2b431392-f2f5-4ead-aca5-bdecef3210cb.md"
$wsOverview.Range("B6").Value = $statusText
$wsOverview.Range("C6").Value = $statusText
$wsOverview.Hyperlinks.Add($wsOverview.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/ffedf167213f298276d27b818585efe18c2bf979/e2e/2b431392-f2f5-4ead-aca5-bdecef3210cb.md", "", "", "2b431392-f2f5-4ead-aca5-bdecef3210cb.md")
$wsOverview.Range("A6").Style = "HyperLink"

$wsOverview.Range("A7").Value = "310c7c43-7abb-4f88-8dac-82b8c4bfe446.md"
$wsOverview.Range("B7").Value = $statusText
$wsOverview.Range("C7").Value = $statusText
$wsOverview.Hyperlinks.Add($wsOverview.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/ffedf167213f298276d27b818585efe18c2bf979/e2e/310c7c43-7abb-4f88-8dac-82b8c4bfe446.md", "", "", "310c7c43-7abb-4f88-8dac-82b8c4bfe446.md")
$wsOverview.Range("A7").Style = "HyperLink"

# ---------------------------------------------------------------------------
# Sheet "zh-cn": two new rows (6 and 7)
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Row 6 - 2b431392-...
$wsZh.Range("A6").Value = "2b431392-f2f5-4ead-aca5-bdecef3210cb.md"
$wsZh.Range("B6").Value = $statusText
$wsZh.Range("C6").Value = "2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.zh-cn.xlf"
$wsZh.Range("D6").Value = "2016-01-28 05:59:50"
$wsZh.Range("E6").Value = "2b431392-f2f5-4ead-aca5-bdecef3210cb.md"
$wsZh.Range("F6").Value = "2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.zh-cn.xlf"
$wsZh.Range("G6").Value = "2016-01-28 06:00:34"
$wsZh.Range("H6").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/ffedf167213f298276d27b818585efe18c2bf979/e2e/2b431392-f2f5-4ead-aca5-bdecef3210cb.md", "", "", "2b431392-f2f5-4ead-aca5-bdecef3210cb.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b73bc382c3f9b832b82cac15f26298cfcaba2b92/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.zh-cn.xlf", "", "", "2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c219b6949431986e4334694d9e3b7be55fa31174/e2e/2b431392-f2f5-4ead-aca5-bdecef3210cb.md", "", "", "2b431392-f2f5-4ead-aca5-bdecef3210cb.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bef468f3a3c6667cddd5cd19539b27c58c74eebe/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.zh-cn.xlf", "", "", "2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.zh-cn.xlf")

$wsZh.Range("A6").Style = "HyperLink"
$wsZh.Range("C6").Style = "HyperLink"
$wsZh.Range("E6").Style = "HyperLink"
$wsZh.Range("F6").Style = "HyperLink"
$wsZh.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 7 - 310c7c43-...
$wsZh.Range("A7").Value = "310c7c43-7abb-4f88-8dac-82b8c4bfe446.md"
$wsZh.Range("B7").Value = $statusText
$wsZh.Range("C7").Value = "310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.zh-cn.xlf"
$wsZh.Range("D7").Value = "2016-01-28 05:59:50"
$wsZh.Range("E7").Value = "310c7c43-7abb-4f88-8dac-82b8c4bfe446.md"
$wsZh.Range("F7").Value = "310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.zh-cn.xlf"
$wsZh.Range("G7").Value = "2016-01-28 06:00:34"
$wsZh.Range("H7").Value = "Include"

$wsZh.Hyperlinks.Add($wsZh.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/ffedf167213f298276d27b818585efe18c2bf979/e2e/310c7c43-7abb-4f88-8dac-82b8c4bfe446.md", "", "", "310c7c43-7abb-4f88-8dac-82b8c4bfe446.md")
$wsZh.Hyperlinks.Add($wsZh.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b73bc382c3f9b832b82cac15f26298cfcaba2b92/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.zh-cn.xlf", "", "", "310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.zh-cn.xlf")
$wsZh.Hyperlinks.Add($wsZh.Range("E7"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/c219b6949431986e4334694d9e3b7be55fa31174/e2e/310c7c43-7abb-4f88-8dac-82b8c4bfe446.md", "", "", "310c7c43-7abb-4f88-8dac-82b8c4bfe446.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/bef468f3a3c6667cddd5cd19539b27c58c74eebe/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/310c7c43-7abb-4f88-dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.zh-cn.xlf", "", "", "310c7c43-7abb-4f88-dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.zh-cn.xlf")

$wsZh.Range("A7").Style = "HyperLink"
$wsZh.Range("C7").Style = "HyperLink"
$wsZh.Range("E7").Style = "HyperLink"
$wsZh.Range("F7").Style = "HyperLink"
$wsZh.Range("D7").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "de-de": two new rows (6 and 7)
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Row 6 - 2b431392-...
$wsDe.Range("A6").Value = "2b431392-f2f5-4ead-aca5-bdecef3210cb.md"
$wsDe.Range("B6").Value = $statusText
$wsDe.Range("C6").Value = "2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.de-de.xlf"
$wsDe.Range("D6").Value = "2016-01-28 06:00:02"
$wsDe.Range("E6").Value = "2b431392-f2f5-4ead-aca5-bdecef3210cb.md"
$wsDe.Range("F6").Value = "2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.de-de.xlf"
$wsDe.Range("G6").Value = "2016-01-28 06:00:53"
$wsDe.Range("H6").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/ffedf167213f298276d27b818585efe18c2bf979/e2e/2b431392-f2f5-4ead-aca5-bdecef3210cb.md", "", "", "2b431392-f2f5-4ead-aca5-bdecef3210cb.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C6"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25400aeee798a320ea462dfbcc625c51a5a62fba/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.de-de.xlf", "", "", "2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("E6"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c6de5ebfbcb58a70fcad64c3cd826c5b48b7fe18/e2e/2b431392-f2f5-4ead-aca5-bdecef3210cb.md", "", "", "2b431392-f2f5-4ead-aca5-bdecef3210cb.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F6"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c780e5905a33d896a5dd3035f500afc68da8657e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.de-de.xlf", "", "", "2b431392-f2f5-4ead-aca5-bdecef3210cb.0e1bf8b158c5d9ae3886aad46243e6d6d22c98ba.de-de.xlf")

$wsDe.Range("A6").Style = "HyperLink"
$wsDe.Range("C6").Style = "HyperLink"
$wsDe.Range("E6").Style = "HyperLink"
$wsDe.Range("F6").Style = "HyperLink"
$wsDe.Range("D6").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 7 - 310c7c43-...
$wsDe.Range("A7").Value = "310c7c43-7abb-4f88-dac-82b8c4bfe446.md"
$wsDe.Range("A7").Value = "310c7c43-7abb-4f88-8dac-82b8c4bfe446.md"
$wsDe.Range("B7").Value = $statusText
$wsDe.Range("C7").Value = "310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.de-de.xlf"
$wsDe.Range("D7").Value = "2016-01-28 06:00:02"
$wsDe.Range("E7").Value = "310c7c43-7abb-4f88-8dac-82b8c4bfe446.md"
$wsDe.Range("F7").Value = "310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.de-de.xlf"
$wsDe.Range("G7").Value = "2016-01-28 06:00:53"
$wsDe.Range("H7").Value = "Include"

$wsDe.Hyperlinks.Add($wsDe.Range("A7"), "https://github.com/OpenLocalizationTest/oltest/blob/ffedf167213f298276d27b818585efe18c2bf979/e2e/310c7c43-7abb-4f88-8dac-82b8c4bfe446.md", "", "", "310c7c43-7abb-4f88-8dac-82b8c4bfe446.md")
$wsDe.Hyperlinks.Add($wsDe.Range("C7"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/25400aeee798a320ea462dfbcc625c51a5a62fba/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.de-de.xlf", "", "", "310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.de-de.xlf")
$wsDe.Hyperlinks.Add($wsDe.Range("E7"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c6de5ebfbcb58a70fcad64c3cd826c5b48b7fe18/e2e/310c7c43-7abb-4f88-8dac-82b8c4bfe446.md", "", "", "310c7c43-7abb-4f88-8dac-82b8c4bfe446.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F7"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/c780e5905a33d896a5dd3035f500afc68da8657e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.de-de.xlf", "", "", "310c7c43-7abb-4f88-8dac-82b8c4bfe446.5af00b2a1e6abce6e4ece677d7a1bc457b1a15ed.de-de.xlf")

$wsDe.Range("A7").Style = "HyperLink"
$wsDe.Range("C7").Style = "HyperLink"
$wsDe.Range("E7").Style = "HyperLink"
$wsDe.Range("F7").Style = "HyperLink"
$wsDe.Range("D7").NumberFormat = "yyyy-mm-dd HH:mm:ss"

Write-Host "Report generated for handback: added 2b431392-f2f5-4ead-aca5-bdecef3210cb and 310c7c43-7abb-4f88-8dac-82b8c4bfe446"
